$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = '42.914.11'
$ws.Cells.Item(2, 5).Value = '  -0.01%  '

# Row 3
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = '2.359.22'
$ws.Cells.Item(3, 5).Value = '  +1.52%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  +0.10%  '

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '302.23'
$ws.Cells.Item(5, 5).Value = '  +0.26%  '

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '95.44'
$ws.Cells.Item(6, 5).Value = '  -0.33%  '

# Row 7
$ws.Cells.Item(7, 5).Value = '  -0.02%  '

# Row 8
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.500'
$ws.Cells.Item(8, 5).Value = '  -0.56%  '

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.481'
$ws.Cells.Item(9, 5).Value = '  -2.45%  '

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '33.96'
$ws.Cells.Item(10, 5).Value = '  -0.98%  '

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.125'
$ws.Cells.Item(11, 5).Value = '  +3.70%  '

# Row 12
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '0.0783'
$ws.Cells.Item(12, 5).Value = '  -0.27%  '

# Row 13
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '18.49'
$ws.Cells.Item(13, 5).Value = '  -2.72%  '

# Row 14
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '2.728.53'
$ws.Cells.Item(14, 5).Value = '  +1.55%  '

# Row 15
$ws.Cells.Item(15, 5).Value = '  -0.88%  '

# Row 16
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '2.359.23'
$ws.Cells.Item(16, 5).Value = '  +1.86%  '

# Row 17
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '0.788'
$ws.Cells.Item(17, 5).Value = '  -0.55%  '

# Row 18
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '42.880.51'
$ws.Cells.Item(18, 5).Value = '  +0.05%  '

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '11.91'
$ws.Cells.Item(19, 5).Value = '  -2.91%  '

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '6.26'
$ws.Cells.Item(20, 5).Value = '  +1.79%  '

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '0.0₃0883'
$ws.Cells.Item(21, 5).Value = '  -0.87%  '

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '67.93'
$ws.Cells.Item(22, 5).Value = '  -0.03%  '

# Row 23
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '234.91'
$ws.Cells.Item(23, 5).Value = '  -0.37%  '

# Row 24
$ws.Cells.Item(24, 5).Value = '  -4.15%  '

# Row 25
$ws.Cells.Item(25, 5).Value = '  +0.05%  '

# Row 26
$ws.Cells.Item(26, 5).Value = '  +0.30%  '

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '24.40'
$ws.Cells.Item(27, 5).Value = '  -0.14%  '

# Row 28
$ws.Cells.Item(28, 5).Value = '  +0.78%  '

# Row 29
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '9.32'
$ws.Cells.Item(29, 5).Value = '  +2.04%  '

# Row 30
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '31.92'
$ws.Cells.Item(30, 5).Value = '  -1.14%  '

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '0.999'
$ws.Cells.Item(31, 5).Value = '  -0.06%  '

# Row 32
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '5.00'
$ws.Cells.Item(32, 5).Value = '  -0.27%  '

# Row 33
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '17.48'
$ws.Cells.Item(33, 5).Value = '  -1.64%  '

# Row 34
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '0.0718'
$ws.Cells.Item(34, 5).Value = '  +2.51%  '

# Row 35
$ws.Cells.Item(35, 2).Value = 'Monero'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '127.11'
$ws.Cells.Item(35, 5).Value = '  -13.19%  '

# Row 36
$ws.Cells.Item(36, 2).Value = 'ARBITRUM'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '1.83'
$ws.Cells.Item(36, 5).Value = '  +1.28%  '

# Row 37
$ws.Cells.Item(37, 5).Value = '  +3.02%  '

# Row 38
$ws.Cells.Item(38, 5).Value = '  -2.84%  '

# Row 39
$ws.Cells.Item(39, 2).Value = 'LidoDAOToken'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '2.81'
$ws.Cells.Item(39, 5).Value = '  +2.54%  '

# Row 40
$ws.Cells.Item(40, 2).Value = 'WEMIXToken'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '2.27'
$ws.Cells.Item(40, 5).Value = '  -2.27%  '

# Row 41
$ws.Cells.Item(41, 5).Value = '  -0.98%  '

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '21.13'
$ws.Cells.Item(42, 5).Value = '  -4.00%  '

# Row 43
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '1.928.23'
$ws.Cells.Item(43, 5).Value = '  +0.13%  '

# Row 44
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '0.0277'
$ws.Cells.Item(44, 5).Value = '  -0.34%  '

# Row 45
$ws.Cells.Item(45, 2).Value = 'FraxShare'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '9.18'
$ws.Cells.Item(45, 5).Value = '  -9.49%  '

# Row 46
$ws.Cells.Item(46, 2).Value = 'NEARProtocol'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '2.69'
$ws.Cells.Item(46, 5).Value = '  -2.22%  '

# Row 47
$ws.Cells.Item(47, 2).Value = 'RocketPoolETH'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '2.589.99'
$ws.Cells.Item(47, 5).Value = '  +1.37%  '

# Row 48
$ws.Cells.Item(48, 2).Value = 'Stacks'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '1.50'
$ws.Cells.Item(48, 5).Value = '  +1.35%  '

# Row 49
$ws.Cells.Item(49, 2).Value = 'TrustWalletToken'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '1.14'
$ws.Cells.Item(49, 5).Value = '  +1.20%  '

# Row 50
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '71.40'
$ws.Cells.Item(50, 5).Value = '  -1.89%  '

# Row 51
$ws.Cells.Item(51, 2).Value = 'MultiversX'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '51.33'
$ws.Cells.Item(51, 5).Value = '  -4.20%  '
